# Update database: drop the oldest period column (6 ماهه منتهی به 1399/06)
# and append the newest period column (12 ماهه منتهی به 1401/12) with its
# figures, shifting every other quarter one column to the left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the oldest data column (D). Everything to its right (E:M) shifts
#    left by one column, taking its values/styles along with it.
$ws.Columns("D:D").Delete()

# 2) The new rightmost column (M) is now blank; clone the formatting
#    (fill/border/alignment/font) from the column immediately to its left (L)
#    so the new column matches the look of the rest of the table.
$ws.Range("L1:L28").Copy()
$ws.Range("M1:M28").PasteSpecial(-4122)
$ws.Columns("M:M").ColumnWidth = 28.17

# 3) The label that used to sit in J9 ("1401-10-28 (6)") shifted left to I9
#    during the delete; the newest report's re-published date replaces it.
$ws.Range("I9").Value = "1402-02-28 (7)"

# 4) Fill in the headers for the newly appended period column.
$ws.Range("M8").Value = "12 ماهه منتهی به 1401/12"
$ws.Range("M9").Value = "1402-02-28"

# 5) Fill in the financial figures for the newly appended period column.
$ws.Range("M11").Value = 23131160
$ws.Range("M12").Value = -13518131
$ws.Range("M13").Value = 9613029
$ws.Range("M14").Value = -339554
$ws.Range("M15").Value = 0
$ws.Range("M16").Value = -154184
$ws.Range("M17").Value = 9119291
$ws.Range("M18").Value = -1937565
$ws.Range("M19").Value = 837932
$ws.Range("M20").Value = 8019658
$ws.Range("M21").Value = -899914
$ws.Range("M22").Value = 7119744
$ws.Range("M23").Value = 0
$ws.Range("M24").Value = 7119744
$ws.Range("M25").Value = 1082
$ws.Range("M26").Value = 6580000
$ws.Range("M27").Value = 1082
